$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the automation-result columns for row 2 (Estado / Transaccion / Fecha)
$ws.Range("E2").Value = "FAILED"

# "Transaccion" came back blank - write an empty text value (quote-prefixed empty
# string) rather than Value="" so a real (empty) text cell is created instead of
# being cleared, then drop back to the default Normal style.
$ws.Range("F2").Value = "'"
$ws.Range("F2").Style = "Normal"

$ws.Range("G2").Value = "14 jul. 2023, 09:06:04"

# Leave the selection where the user ended up after running the test
$ws.Range("A2").Select()
